$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 11737.5
$ws.Range("J48").Value = 3475
$ws.Range("L48").Value = 10425
$ws.Range("N48").Value = -11009

$ws.Range("H56").Value = 11737.5
$ws.Range("J56").Value = 3475
$ws.Range("L56").Value = 10425
$ws.Range("N56").Value = -11493

$ws.Range("H112").Value = 2745.1714
$ws.Range("I112").Value = 2098.5
$ws.Range("J112").Value = 2828.6128
$ws.Range("K112").Value = 6295.5
$ws.Range("L112").Value = 8485.838400000001
$ws.Range("M112").Value = -5187.5
$ws.Range("N112").Value = -10701.8384

$ws.Range("H127").Value = 2732.45
$ws.Range("I127").Value = 2732.45
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 8197.349999999999
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -3237.349999999999
$ws.Range("N127").ClearContents()

$ws.Range("H135").Value = 2961.5
$ws.Range("I135").Value = 2053.8
$ws.Range("J135").Value = 7500
$ws.Range("K135").Value = 18484.2
$ws.Range("L135").Value = 67500
$ws.Range("M135").Value = -15949.2
$ws.Range("N135").Value = -72570

$ws.Range("H137").Value = 4600.4517
$ws.Range("I137").Value = 7338.1177
$ws.Range("J137").Value = 1276.1428
$ws.Range("K137").Value = 22014.3531
$ws.Range("L137").Value = 3828.4284
$ws.Range("M137").Value = -19464.3531
$ws.Range("N137").Value = -8928.428400000001

$ws.Range("H138").Value = 2871.6223
$ws.Range("I138").Value = 3266.1538
$ws.Range("J138").Value = 2711.3438
$ws.Range("K138").Value = 9798.4614
$ws.Range("L138").Value = 8134.0314
$ws.Range("M138").Value = -4658.4614
$ws.Range("N138").Value = -18414.0314

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9363.647999999999
$ws.Range("I32").Value = 8861.582
$ws.Range("J32").Value = 14169.143
$ws.Range("K32").Value = 8861.582
$ws.Range("L32").Value = 14169.143
$ws.Range("M32").Value = -8574.582
$ws.Range("N32").Value = -14743.143

$ws.Range("H61").Value = 6611.4
$ws.Range("J61").Value = 5106.75
$ws.Range("L61").Value = 5106.75
$ws.Range("N61").Value = -5530.75

$ws.Range("H74").Value = 3745.1667
$ws.Range("I74").Value = 1942.25
$ws.Range("J74").Value = 7351
$ws.Range("K74").Value = 1942.25
$ws.Range("L74").Value = 7351
$ws.Range("M74").Value = -1068.25
$ws.Range("N74").Value = -9099

$ws.Range("H77").Value = 3745.1667
$ws.Range("I77").Value = 1942.25
$ws.Range("J77").Value = 7351
$ws.Range("K77").Value = 9711.25
$ws.Range("L77").Value = 36755
$ws.Range("M77").Value = -5343.25
$ws.Range("N77").Value = -45491

$ws.Range("H119").Value = 399999
$ws.Range("J119").Value = 399999
$ws.Range("L119").Value = 399999
$ws.Range("N119").Value = -409675

$ws.Range("H132").Value = 3026.2808
$ws.Range("I132").Value = 3012.7073
$ws.Range("K132").Value = 9038.1219
$ws.Range("M132").Value = -6508.1219

$ws.Range("H136").Value = 6611.4
$ws.Range("J136").Value = 5106.75
$ws.Range("L136").Value = 15320.25
$ws.Range("N136").Value = -20420.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 13975
$ws.Range("I26").Value = 13975
$ws.Range("K26").Value = 13975
$ws.Range("M26").Value = -13683

$ws.Range("H128").Value = 4999.375
$ws.Range("I128").Value = 4999.375
$ws.Range("K128").Value = 14998.125
$ws.Range("M128").Value = -12508.125

$ws.Range("H134").Value = 3460.3948
$ws.Range("I134").Value = 2750.375
$ws.Range("K134").Value = 8251.125
$ws.Range("M134").Value = -5716.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4209
$ws.Range("J31").Value = 10411.111
$ws.Range("L31").Value = 10411.111
$ws.Range("N31").Value = -11001.111

$ws.Range("H34").Value = 4209
$ws.Range("J34").Value = 10411.111
$ws.Range("L34").Value = 10411.111
$ws.Range("N34").Value = -10815.111

$ws.Range("H58").Value = 1436.6857
$ws.Range("J58").Value = 2859.625
$ws.Range("L58").Value = 2859.625
$ws.Range("N58").Value = -3265.625

$ws.Range("H94").Value = 7070.4
$ws.Range("I94").Value = 7984.8887
$ws.Range("J94").Value = 5698.6665
$ws.Range("K94").Value = 7984.8887
$ws.Range("L94").Value = 5698.6665
$ws.Range("M94").Value = -7533.8887
$ws.Range("N94").Value = -6600.6665

$ws.Range("H107").Value = 1066.5333
$ws.Range("I107").Value = 851.5
$ws.Range("J107").Value = 1254.6875
$ws.Range("K107").Value = 851.5
$ws.Range("L107").Value = 1254.6875
$ws.Range("M107").Value = 1068.5
$ws.Range("N107").Value = -5094.6875

$ws.Range("H121").Value = 64325
$ws.Range("J121").Value = 64325
$ws.Range("L121").Value = 64325
$ws.Range("N121").Value = -66945

$ws.Range("H132").Value = 5457.4443
$ws.Range("I132").Value = 1375.6842
$ws.Range("J132").Value = 15151.625
$ws.Range("K132").Value = 4127.0526
$ws.Range("L132").Value = 45454.875
$ws.Range("M132").Value = -1597.0526
$ws.Range("N132").Value = -50514.875

$ws.Range("H136").Value = 1436.6857
$ws.Range("J136").Value = 2859.625
$ws.Range("L136").Value = 8578.875
$ws.Range("N136").Value = -13678.875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 99.36
$ws.Range("J2").Value = 134.6
$ws.Range("L2").Value = 807.5999999999999
$ws.Range("N2").Value = -1033.6

$ws.Range("H3").Value = 7659
$ws.Range("I3").Value = 5122.9165
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 15368.7495
$ws.Range("L3").Value = 30000
$ws.Range("M3").Value = -15256.7495
$ws.Range("N3").Value = -30224

$ws.Range("H23").Value = 243.86667
$ws.Range("I23").Value = 251.6
$ws.Range("J23").Value = 240
$ws.Range("K23").Value = 754.8
$ws.Range("L23").Value = 720
$ws.Range("M23").Value = -519.8
$ws.Range("N23").Value = -1190

$ws.Range("H103").Value = 1378.8182
$ws.Range("I103").Value = 334.25
$ws.Range("K103").Value = 1002.75
$ws.Range("M103").Value = -123.75

$ws.Range("H113").Value = 463.6875
$ws.Range("J113").Value = 420.2
$ws.Range("L113").Value = 1260.6
$ws.Range("N113").Value = -5600.6

$ws.Range("H117").Value = 3771.889
$ws.Range("J117").Value = 2999.5
$ws.Range("L117").Value = 8998.5
$ws.Range("N117").Value = -15882.5

$ws.Range("H121").Value = 1165.0952
$ws.Range("J121").Value = 1418.6666
$ws.Range("L121").Value = 4255.9998
$ws.Range("N121").Value = -6875.9998

$ws.Range("H131").Value = 1401447.9
$ws.Range("I131").Value = 1839025.4
$ws.Range("J131").Value = 1199.8
$ws.Range("K131").Value = 5517076.199999999
$ws.Range("L131").Value = 3599.4
$ws.Range("M131").Value = -5512036.199999999
$ws.Range("N131").Value = -13679.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H120").Value = 160877.67
$ws.Range("J120").Value = 160877.67
$ws.Range("L120").Value = 160877.67
$ws.Range("N120").Value = -170553.67

$ws.Range("H132").Value = 2295.0625
$ws.Range("I132").Value = 1848.3334
$ws.Range("K132").Value = 5545.0002
$ws.Range("M132").Value = -3015.0002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1744.4688
$ws.Range("I40").Value = 1484.3462
$ws.Range("J40").Value = 2871.6667
$ws.Range("K40").Value = 1484.3462
$ws.Range("L40").Value = 2871.6667
$ws.Range("M40").Value = -1348.3462
$ws.Range("N40").Value = -3143.6667

$ws.Range("H132").Value = 10373.158
$ws.Range("I132").Value = 6882.091
$ws.Range("K132").Value = 20646.273
$ws.Range("M132").Value = -18116.273

$ws.Range("H134").Value = 98089.8
$ws.Range("J134").Value = 97612.25
$ws.Range("L134").Value = 97612.25
$ws.Range("N134").Value = -107752.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 483.35483
$ws.Range("I113").Value = 414.9
$ws.Range("J113").Value = 607.8182
$ws.Range("K113").Value = 1244.7
$ws.Range("L113").Value = 1823.4546
$ws.Range("M113").Value = 925.3000000000002
$ws.Range("N113").Value = -6163.4546

$ws.Range("H122").Value = 6633.354
$ws.Range("I122").Value = 2456.244
$ws.Range("K122").Value = 7368.732
$ws.Range("M122").Value = -4918.732

$ws.Range("H132").Value = 3279.4348
$ws.Range("I132").Value = 3087.2058
$ws.Range("J132").Value = 3824.0833
$ws.Range("K132").Value = 9261.617400000001
$ws.Range("L132").Value = 11472.2499
$ws.Range("M132").Value = -6731.617400000001
$ws.Range("N132").Value = -16532.2499

$ws.Range("H136").Value = 1993.9783
$ws.Range("I136").Value = 1588.641
$ws.Range("K136").Value = 4765.923000000001
$ws.Range("M136").Value = -2215.923000000001

$ws.Range("H137").Value = 144416.17
$ws.Range("J137").Value = 144416.17
$ws.Range("L137").Value = 144416.17
$ws.Range("N137").Value = -154616.17
